$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'299.16"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-1.80%"
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'31.43"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-1.56%"
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").Value = "'5.094"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-1.80%"
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'0.07937"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'6.14%"
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'2.271"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-5.31%"
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").Value = "'7.740"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-3.33%"
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'3.865"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'0.08%"
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").Value = "'0.9169"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-0.07%"
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").Value = "'0.1732"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-0.12%"
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").Value = "'0.07347"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-4.15%"
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'0.09073"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'11.15%"
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").Value = "'0.03025"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'0.40%"
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'0.1002"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'0.82%"
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").Value = "'0.001516"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'0.64%"
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'0.006059"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-0.01%"
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("E17").Value = "'-0.79%"
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").Value = "'2.265"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'1.66%"
$ws.Range("E18").ClearFormats()

# Row 20
$ws.Range("E20").Value = "'-2.12%"
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").Value = "'4.194"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'-9.85%"
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").Value = "'0.1699"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'8.67%"
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("E23").Value = "'0.32%"
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").Value = "'0.001240"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-1.51%"
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").Value = "'0.004459"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-1.55%"
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("E26").Value = "'-7.60%"
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").Value = "'24.03%"
$ws.Range("E27").ClearFormats()

# Row 39
$ws.Range("D39").Value = "'0.01746"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'-1.74%"
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").Value = "'0.04607"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'0.89%"
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("E41").Value = "'-6.25%"
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").Value = "'0.1356"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-0.47%"
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").Value = "'0.002188"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'0.10%"
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").Value = "'0.009539"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-11.63%"
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").Value = "'0.00006269"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-2.58%"
$ws.Range("E45").ClearFormats()

# Row 47
$ws.Range("E47").Value = "'-19.34%"
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").Value = "'1.159"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'41.27%"
$ws.Range("E48").ClearFormats()
